# The commit inserts one new observation row ("semanal" data point) into the
# Cilantro price history sheet. It is inserted right before the existing
# row 37, which shifts every subsequent row (old 37..139) down by one
# (new 38..140), growing the used range from A1:R139 to A1:R140.
#
# The newly inserted row 37 carries the constant attribute columns shared by
# every other row in the table, a new date (serial 44672 = 2022-04-21), an
# unchanged Volumen (3000) and the min/max/avg/unit-price columns bumped to
# 2000 / 2500 / 2250 / 1500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 37:139 down to 38:140, leaving a blank row 37 (inheriting the
# formatting of the row above, same as Excel's native "Insert Row").
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the new data point.
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C37").Value = 'Coquimbo'
$ws.Range("D37").Value = 44672
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100112040
$ws.Range("G37").Value = 'Cilantro'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 2500
$ws.Range("M37").Value = 2250
$ws.Range("N37").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O37").Value = 'Provincia del Elquí'
$ws.Range("P37").Value = 1500
$ws.Range("Q37").Value = 1.5
$ws.Range("R37").Value = 'Hortaliza'
